$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "Name"
$ws.Range("G4").Value = "S101"
$ws.Range("G7").Value = "Q101"
$ws.Range("G11").Value = "LED101"
$ws.Range("G16").Value = "IC101"
$ws.Range("G32").Value = "C303"
$ws.Range("G23").Value = "C202"
$ws.Range("G37").Value = "R202"
$ws.Range("G5").Value = "X202"
$ws.Range("G6").Value = "X201"
$ws.Range("G8").Value = "J101"
$ws.Range("G9").Value = "J401, J402"
$ws.Range("G10").Value = "J403, J404"
$ws.Range("G12").Value = "LED102"
$ws.Range("G13").Value = "LED201, LED301"
$ws.Range("G14").Value = "LED202, LED302"
$ws.Range("G15").Value = "IC203"
$ws.Range("G17").Value = "IC202"
$ws.Range("G18").Value = "IC201"
$ws.Range("G19").Value = "MOD301"
$ws.Range("G20").Value = "D201, D202, D203, D204"
$ws.Range("G21").Value = "L202, L301, L302"
$ws.Range("G22").Value = "L201"
$ws.Range("G24").Value = "C203"
$ws.Range("G25").Value = "C204, C209, C215, C216, C304"
$ws.Range("G26").Value = "C109"
$ws.Range("G27").Value = "C104, C106, C301"
$ws.Range("G28").Value = "C103, C105, C107, C108, C110, C111, C201, C205, C207, C210, C213, C214, C217, C302, C305"
$ws.Range("G29").Value = "C206, C208"
$ws.Range("G30").Value = "C211, C212"
$ws.Range("G31").Value = "C101, C102"
$ws.Range("G33").Value = "R101, R102, R104, R105"
$ws.Range("G34").Value = "R205, R206"
$ws.Range("G35").Value = "R107, R108, R204, R207, R208, R301, R302"
$ws.Range("G36").Value = "R103, R203"
$ws.Range("G38").Value = "R201"

$ws.Range("G4:G38").NumberFormat = "@"

$ws.Columns("G").ColumnWidth = 78.14

$ws.Range("G5").Select() | Out-Null
